$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.5427046263345195
$summary.Range("C2").Value = 0.08960573476702509
$summary.Range("D2").Value = 0.8928571428571429
$summary.Range("E2").Value = 0.1628664495114006
$summary.Range("F2").Value = 0.319693094629156
$summary.Range("G2").Value = 0.6639427987742594
$summary.Range("H2").Value = 0.690342429106474
$summary.Range("I2").Value = 25
$summary.Range("J2").Value = 254
$summary.Range("K2").Value = 280
$summary.Range("L2").Value = 3

# --- Sheet 2: Classification Report ---
$report = $wb.Worksheets.Item("Classification Report")

# Row 2 (class "0")
$report.Range("B2").Value = 0.9893992932862191
$report.Range("C2").Value = 0.5243445692883895
$report.Range("D2").Value = 0.6854345165238678

# Row 3 (class "1")
$report.Range("B3").Value = 0.08960573476702509
$report.Range("C3").Value = 0.8928571428571429
$report.Range("D3").Value = 0.1628664495114006

# Row 4 (accuracy)
$report.Range("B4").Value = 0.5427046263345195
$report.Range("C4").Value = 0.5427046263345195
$report.Range("D4").Value = 0.5427046263345195
$report.Range("E4").Value = 0.5427046263345195

# Row 5 (macro avg)
$report.Range("B5").Value = 0.5395025140266221
$report.Range("C5").Value = 0.7086008560727661
$report.Range("D5").Value = 0.4241504830176342

# Row 6 (weighted avg)
$report.Range("B6").Value = 0.9445697209756542
$report.Range("C6").Value = 0.5427046263345195
$report.Range("D6").Value = 0.6593990968150616

# --- Sheet 3: Confusion Matrix ---
$confusion = $wb.Worksheets.Item("Confusion Matrix")
$confusion.Range("B2").Value = 280
$confusion.Range("C2").Value = 254
$confusion.Range("B3").Value = 3
$confusion.Range("C3").Value = 25
